$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet lists each data extract along with the folder path(s) it has
# historically lived in. Row 3 is a duplicate WHO extract entry (same
# nm/sql_str as row 4, just pointing at an old, no-longer-used network
# path) - remove it.

# Preserve the hyperlink that lives on the "modeling" row further down
# (currently D11) so it survives the row shift caused by deleting row 3.
$hyperlinkTarget = $ws.Range("D11").Hyperlinks.Item(1).Address
if ([string]::IsNullOrEmpty($hyperlinkTarget)) {
    $hyperlinkTarget = "https://storphacidpcbns02.blob.core.windows.net/hcdaily"
}
$ws.Range("D11").Hyperlinks.Delete()

# Delete the whole duplicate row; everything below shifts up one row.
$ws.Rows.Item(3).Delete()

# Re-create the hyperlink at its new location (was D11, now D10).
$ws.Hyperlinks.Add($ws.Range("D10"), $hyperlinkTarget)
# Adding a hyperlink auto-applies the built-in blue/underline "Hyperlink"
# look; the cell kept its plain formatting in the original file, so
# restore it.
$ws.Range("D10").Style = "Normal"

# Select the row that now occupies row 3 (mirrors selecting/deleting a
# whole row in the Excel UI, which leaves the shifted-up row selected).
$ws.Rows.Item(3).Select() | Out-Null
